# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46062
$ws.Range("B2").Value = 0.48
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = -0.02
$ws.Range("F2").Value = -0.02
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.02
$ws.Range("I2").Value = 0.83
$ws.Range("J2").Value = 5.09
$ws.Range("K2").Value = 19.84
$ws.Range("L2").Value = 12.32
$ws.Range("M2").Value = 3.22
$ws.Range("N2").Value = 0.71
$ws.Range("O2").Value = 0.1
$ws.Range("P2").Value = 0.02
$ws.Range("Q2").Value = 0.08
$ws.Range("R2").Value = 0.97
$ws.Range("S2").Value = 3.35
$ws.Range("T2").Value = 2.4
$ws.Range("U2").Value = 4.84
$ws.Range("V2").Value = 10.36
$ws.Range("W2").Value = 8.449999999999999
$ws.Range("X2").Value = 3.79
$ws.Range("Y2").Value = 1.09
$ws.Range("Z2").Value = 3.25
$ws.Range("AA2").Value = "8h-12h"
$ws.Range("AB2").Value = 10.12
$ws.Range("AC2").Value = "8h-10h"
$ws.Range("AD2").Value = 12.46
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 9.4
$ws.Range("AG2").Value = "0h-23h"
